$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AF2").ClearContents()
$ws.Range("A2").Value = 17177744
$ws.Range("B2").Value = 93145
$ws.Range("C2").Value = 'Ovaliderad'
$ws.Range("D2").Value = 'LC'
$ws.Range("E2").Value = 2667
$ws.Range("F2").Value = 'Platt fjädermossa'
$ws.Range("G2").Value = 'Neckera complanata'
$ws.Range("H2").Value = '(Hedw.) Huebener'
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("P2").Value = 'S Kramnäs, Srm'
$ws.Range("Q2").Value = 596512.1226172579
$ws.Range("R2").Value = 6550345.023338513
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 'Södermanland'
$ws.Range("U2").Value = 'Flen'
$ws.Range("V2").Value = 'Södermanland'
$ws.Range("W2").Value = 'Mellösa'
$ws.Range("Y2").Value = '2015-03-02'
$ws.Range("Z2").Value = '00:00'
$ws.Range("AA2").Value = '2015-03-02'
$ws.Range("AB2").Value = '00:00'
$ws.Range("AC2").Value = 'På block och lodyta. Nyckelbiotopsstatus.'
$ws.Range("AD2").Value = $False
$ws.Range("AE2").Value = $False
$ws.Range("AG2").Value = $False
$ws.Range("AT2").Value = ""
$ws.Range("AW2").Value = 'Bo Törnquist'
$ws.Range("AX2").Value = 'Bo Törnquist'
$ws.Range("AY2").Value = ""

# Row 3
$ws.Range("A3").Value = 17177743
$ws.Range("B3").Value = 93132
$ws.Range("C3").Value = 'Ovaliderad'
$ws.Range("D3").Value = 'LC'
$ws.Range("E3").Value = 2671
$ws.Range("F3").Value = 'Fällmossa'
$ws.Range("G3").Value = 'Antitrichia curtipendula'
$ws.Range("H3").Value = '(Hedw.) Brid.'
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("P3").Value = 'S Kramnäs, Srm'
$ws.Range("Q3").Value = 596508.9821110814
$ws.Range("R3").Value = 6550347.00089203
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 'Södermanland'
$ws.Range("U3").Value = 'Flen'
$ws.Range("V3").Value = 'Södermanland'
$ws.Range("W3").Value = 'Mellösa'
$ws.Range("Y3").Value = '2015-03-02'
$ws.Range("Z3").Value = '00:00'
$ws.Range("AA3").Value = '2015-03-02'
$ws.Range("AB3").Value = '00:00'
$ws.Range("AC3").Value = 'Riklig förekomst på block, lodytor och enstaka trädbaser.'
$ws.Range("AD3").Value = $False
$ws.Range("AE3").Value = $False
$ws.Range("AG3").Value = $False
$ws.Range("AT3").Value = ""
$ws.Range("AW3").Value = 'Bo Törnquist'
$ws.Range("AX3").Value = 'Bo Törnquist'
$ws.Range("AY3").Value = ""

# Row 4
$ws.Range("A4").Value = 17177745
$ws.Range("B4").Value = 92939
$ws.Range("C4").Value = 'Ovaliderad'
$ws.Range("D4").Value = 'LC'
$ws.Range("E4").Value = 2779
$ws.Range("F4").Value = 'Guldlockmossa'
$ws.Range("G4").Value = 'Homalothecium sericeum'
$ws.Range("H4").Value = '(Hedw.) Schimp.'
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("P4").Value = 'S Kramnäs, Srm'
$ws.Range("Q4").Value = 596512.1226172579
$ws.Range("R4").Value = 6550345.023338513
$ws.Range("S4").Value = 5
$ws.Range("T4").Value = 'Södermanland'
$ws.Range("U4").Value = 'Flen'
$ws.Range("V4").Value = 'Södermanland'
$ws.Range("W4").Value = 'Mellösa'
$ws.Range("Y4").Value = '2015-03-02'
$ws.Range("Z4").Value = '00:00'
$ws.Range("AA4").Value = '2015-03-02'
$ws.Range("AB4").Value = '00:00'
$ws.Range("AC4").Value = 'Lodyta'
$ws.Range("AD4").Value = $False
$ws.Range("AE4").Value = $False
$ws.Range("AG4").Value = $False
$ws.Range("AT4").Value = ""
$ws.Range("AW4").Value = 'Bo Törnquist'
$ws.Range("AX4").Value = 'Bo Törnquist'
$ws.Range("AY4").Value = ""

# Row 5
$ws.Range("A5").Value = 17177741
$ws.Range("B5").Value = 93132
$ws.Range("C5").Value = 'Ovaliderad'
$ws.Range("D5").Value = 'LC'
$ws.Range("E5").Value = 2671
$ws.Range("F5").Value = 'Fällmossa'
$ws.Range("G5").Value = 'Antitrichia curtipendula'
$ws.Range("H5").Value = '(Hedw.) Brid.'
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("P5").Value = 'S Kramnäs, Srm'
$ws.Range("Q5").Value = 596502.1510387775
$ws.Range("R5").Value = 6550331.917076289
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = 'Södermanland'
$ws.Range("U5").Value = 'Flen'
$ws.Range("V5").Value = 'Södermanland'
$ws.Range("W5").Value = 'Mellösa'
$ws.Range("Y5").Value = '2015-03-02'
$ws.Range("Z5").Value = '00:00'
$ws.Range("AA5").Value = '2015-03-02'
$ws.Range("AB5").Value = '00:00'
$ws.Range("AC5").Value = 'Riklig förekomst på block, lodytor och enstaka trädbaser. Blockrik bergbrant med nyckelbiotopsstatus.'
$ws.Range("AD5").Value = $False
$ws.Range("AE5").Value = $False
$ws.Range("AG5").Value = $False
$ws.Range("AT5").Value = ""
$ws.Range("AW5").Value = 'Bo Törnquist'
$ws.Range("AX5").Value = 'Bo Törnquist'
$ws.Range("AY5").Value = ""

# Row 6
$ws.Range("L6").ClearContents()
$ws.Range("A6").Value = 17177739
$ws.Range("B6").Value = 89940
$ws.Range("C6").Value = 'Ovaliderad'
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 3884
$ws.Range("F6").Value = 'Hasselticka'
$ws.Range("G6").Value = 'Dichomitus campestris'
$ws.Range("H6").Value = '(Quél.) Domański & Orlicz'
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("P6").Value = 'S Kramnäs, Srm'
$ws.Range("Q6").Value = 596515.088978241
$ws.Range("R6").Value = 6550309.1046673
$ws.Range("S6").Value = 5
$ws.Range("T6").Value = 'Södermanland'
$ws.Range("U6").Value = 'Flen'
$ws.Range("V6").Value = 'Södermanland'
$ws.Range("W6").Value = 'Mellösa'
$ws.Range("Y6").Value = '2015-03-02'
$ws.Range("Z6").Value = '00:00'
$ws.Range("AA6").Value = '2015-03-02'
$ws.Range("AB6").Value = '00:00'
$ws.Range("AC6").Value = 'Flera döda stammar.'
$ws.Range("AD6").Value = $False
$ws.Range("AE6").Value = $False
$ws.Range("AF6").Value = ""
$ws.Range("AG6").Value = $False
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = 'Bo Törnquist'
$ws.Range("AX6").Value = 'Bo Törnquist'
$ws.Range("AY6").Value = ""

# Row 7
$ws.Range("A7").Value = 17177742
$ws.Range("B7").Value = 4717
$ws.Range("C7").Value = 'Ovaliderad'
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 102306
$ws.Range("F7").Value = 'Granbarkgnagare'
$ws.Range("G7").Value = 'Microbregma emarginatum'
$ws.Range("H7").Value = '(Duftschmid, 1825)'
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("P7").Value = 'S Kramnäs, Srm'
$ws.Range("Q7").Value = 596494.09243898
$ws.Range("R7").Value = 6550345.08282318
$ws.Range("S7").Value = 5
$ws.Range("T7").Value = 'Södermanland'
$ws.Range("U7").Value = 'Flen'
$ws.Range("V7").Value = 'Södermanland'
$ws.Range("W7").Value = 'Mellösa'
$ws.Range("Y7").Value = '2015-03-02'
$ws.Range("Z7").Value = '00:00'
$ws.Range("AA7").Value = '2015-03-02'
$ws.Range("AB7").Value = '00:00'
$ws.Range("AC7").Value = 'Kläckhål på gammal grov gran.'
$ws.Range("AD7").Value = $False
$ws.Range("AE7").Value = $False
$ws.Range("AG7").Value = $False
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = 'Bo Törnquist'
$ws.Range("AX7").Value = 'Bo Törnquist'
$ws.Range("AY7").Value = ""

# Row 8
$ws.Range("AC8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("A8").Value = 17177746
$ws.Range("B8").Value = 89412
$ws.Range("C8").Value = 'Ovaliderad'
$ws.Range("D8").Value = 'NT'
$ws.Range("E8").Value = 5442
$ws.Range("F8").Value = 'Tallticka'
$ws.Range("G8").Value = 'Porodaedalea pini'
$ws.Range("H8").Value = '(Brot.) Murrill'
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = '2'
$ws.Range("J8").Value = 'fruktkroppar'
$ws.Range("K8").Value = ""
$ws.Range("P8").Value = 'S Kramnäs, Srm'
$ws.Range("Q8").Value = 596514.216958874
$ws.Range("R8").Value = 6550364.101154891
$ws.Range("S8").Value = 5
$ws.Range("T8").Value = 'Södermanland'
$ws.Range("U8").Value = 'Flen'
$ws.Range("V8").Value = 'Södermanland'
$ws.Range("W8").Value = 'Mellösa'
$ws.Range("Y8").Value = '2015-03-02'
$ws.Range("Z8").Value = '00:00'
$ws.Range("AA8").Value = '2015-03-02'
$ws.Range("AB8").Value = '00:00'
$ws.Range("AD8").Value = $False
$ws.Range("AE8").Value = $False
$ws.Range("AF8").Value = ""
$ws.Range("AG8").Value = $False
$ws.Range("AT8").Value = ""
$ws.Range("AW8").Value = 'Bo Törnquist'
$ws.Range("AX8").Value = 'Bo Törnquist'
$ws.Range("AY8").Value = ""

# Row 9
$ws.Range("AF9").ClearContents()
$ws.Range("A9").Value = 17180057
$ws.Range("B9").Value = 93145
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 2667
$ws.Range("F9").Value = 'Platt fjädermossa'
$ws.Range("G9").Value = 'Neckera complanata'
$ws.Range("H9").Value = '(Hedw.) Huebener'
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("P9").Value = 'S Kramnäs, Srm'
$ws.Range("Q9").Value = 596513.8151775175
$ws.Range("R9").Value = 6550380.030807917
$ws.Range("S9").Value = 5
$ws.Range("T9").Value = 'Södermanland'
$ws.Range("U9").Value = 'Flen'
$ws.Range("V9").Value = 'Södermanland'
$ws.Range("W9").Value = 'Mellösa'
$ws.Range("Y9").Value = '2015-03-02'
$ws.Range("Z9").Value = '00:00'
$ws.Range("AA9").Value = '2015-03-02'
$ws.Range("AB9").Value = '00:00'
$ws.Range("AC9").Value = 'På gammal asp och block.'
$ws.Range("AD9").Value = $False
$ws.Range("AE9").Value = $False
$ws.Range("AG9").Value = $False
$ws.Range("AT9").Value = ""
$ws.Range("AW9").Value = 'Bo Törnquist'
$ws.Range("AX9").Value = 'Bo Törnquist'
$ws.Range("AY9").Value = ""

# Row 10
$ws.Range("AF10").ClearContents()
$ws.Range("A10").Value = 17180059
$ws.Range("B10").Value = 93235
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 210
$ws.Range("F10").Value = 'Grön sköldmossa'
$ws.Range("G10").Value = 'Buxbaumia viridis'
$ws.Range("H10").Value = '(Moug. ex Lam. & DC.) Brid. ex Moug. & Nestl.'
$ws.Range("I10").NumberFormat = "@"
$ws.Range("I10").Value = '20'
$ws.Range("J10").Value = 'kapslar'
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("P10").Value = 'S Kramnäs, Srm'
$ws.Range("Q10").Value = 596440.0180564586
$ws.Range("R10").Value = 6550405.936094929
$ws.Range("S10").Value = 5
$ws.Range("T10").Value = 'Södermanland'
$ws.Range("U10").Value = 'Flen'
$ws.Range("V10").Value = 'Södermanland'
$ws.Range("W10").Value = 'Mellösa'
$ws.Range("Y10").Value = '2015-03-02'
$ws.Range("Z10").Value = '00:00'
$ws.Range("AA10").Value = '2015-03-02'
$ws.Range("AB10").Value = '00:00'
$ws.Range("AC10").Value = '20 sporkapslar på gammal granlåga.'
$ws.Range("AD10").Value = $False
$ws.Range("AE10").Value = $False
$ws.Range("AG10").Value = $False
$ws.Range("AT10").Value = ""
$ws.Range("AW10").Value = 'Bo Törnquist'
$ws.Range("AX10").Value = 'Bo Törnquist'
$ws.Range("AY10").Value = ""

# Row 11
$ws.Range("AF11").ClearContents()
$ws.Range("A11").Value = 17180058
$ws.Range("B11").Value = 4717
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 102306
$ws.Range("F11").Value = 'Granbarkgnagare'
$ws.Range("G11").Value = 'Microbregma emarginatum'
$ws.Range("H11").Value = '(Duftschmid, 1825)'
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("P11").Value = 'S Kramnäs, Srm'
$ws.Range("Q11").Value = 596463.7622184914
$ws.Range("R11").Value = 6550403.963688263
$ws.Range("S11").Value = 5
$ws.Range("T11").Value = 'Södermanland'
$ws.Range("U11").Value = 'Flen'
$ws.Range("V11").Value = 'Södermanland'
$ws.Range("W11").Value = 'Mellösa'
$ws.Range("Y11").Value = '2015-03-02'
$ws.Range("Z11").Value = '00:00'
$ws.Range("AA11").Value = '2015-03-02'
$ws.Range("AB11").Value = '00:00'
$ws.Range("AC11").Value = 'Kläckhål på gammal gran.'
$ws.Range("AD11").Value = $False
$ws.Range("AE11").Value = $False
$ws.Range("AG11").Value = $False
$ws.Range("AT11").Value = ""
$ws.Range("AW11").Value = 'Bo Törnquist'
$ws.Range("AX11").Value = 'Bo Törnquist'
$ws.Range("AY11").Value = ""

# Row 12
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("A12").Value = 17162672
$ws.Range("B12").Value = 89376
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 4660
$ws.Range("F12").Value = 'Rävticka'
$ws.Range("G12").Value = 'Inocutis rheades'
$ws.Range("H12").Value = '(Pers.) Fiasson & Niemelä'
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("P12").Value = 'S Kramnäs, Srm'
$ws.Range("Q12").Value = 596445.2567064186
$ws.Range("R12").Value = 6550238.955486102
$ws.Range("S12").Value = 5
$ws.Range("T12").Value = 'Södermanland'
$ws.Range("U12").Value = 'Flen'
$ws.Range("V12").Value = 'Södermanland'
$ws.Range("W12").Value = 'Mellösa'
$ws.Range("Y12").Value = '2015-03-02'
$ws.Range("Z12").Value = '00:00'
$ws.Range("AA12").Value = '2015-03-02'
$ws.Range("AB12").Value = '00:00'
$ws.Range("AC12").Value = 'Klen högstubbe av asp.'
$ws.Range("AD12").Value = $False
$ws.Range("AE12").Value = $False
$ws.Range("AF12").Value = ""
$ws.Range("AG12").Value = $False
$ws.Range("AT12").Value = ""
$ws.Range("AW12").Value = 'Bo Törnquist'
$ws.Range("AX12").Value = 'Bo Törnquist'
$ws.Range("AY12").Value = ""

# Row 13
$ws.Range("AF13").ClearContents()
$ws.Range("A13").Value = 17177735
$ws.Range("B13").Value = 93132
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'LC'
$ws.Range("E13").Value = 2671
$ws.Range("F13").Value = 'Fällmossa'
$ws.Range("G13").Value = 'Antitrichia curtipendula'
$ws.Range("H13").Value = '(Hedw.) Brid.'
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("P13").Value = 'S Kramnäs, Srm'
$ws.Range("Q13").Value = 596484.0673687706
$ws.Range("R13").Value = 6550191.084131638
$ws.Range("S13").Value = 5
$ws.Range("T13").Value = 'Södermanland'
$ws.Range("U13").Value = 'Flen'
$ws.Range("V13").Value = 'Södermanland'
$ws.Range("W13").Value = 'Mellösa'
$ws.Range("Y13").Value = '2015-03-02'
$ws.Range("Z13").Value = '00:00'
$ws.Range("AA13").Value = '2015-03-02'
$ws.Range("AB13").Value = '00:00'
$ws.Range("AC13").Value = 'Riklig förekomst på block, lodytor och enstaka aspstammar.'
$ws.Range("AD13").Value = $False
$ws.Range("AE13").Value = $False
$ws.Range("AG13").Value = $False
$ws.Range("AT13").Value = ""
$ws.Range("AW13").Value = 'Bo Törnquist'
$ws.Range("AX13").Value = 'Bo Törnquist'
$ws.Range("AY13").Value = ""

# Row 14
$ws.Range("A14").Value = 17177734
$ws.Range("B14").Value = 93132
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 2671
$ws.Range("F14").Value = 'Fällmossa'
$ws.Range("G14").Value = 'Antitrichia curtipendula'
$ws.Range("H14").Value = '(Hedw.) Brid.'
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = ""
$ws.Range("L14").Value = ""
$ws.Range("P14").Value = 'S Kramnäs, Srm'
$ws.Range("Q14").Value = 596443.9843037433
$ws.Range("R14").Value = 6550228.125120561
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = 'Södermanland'
$ws.Range("U14").Value = 'Flen'
$ws.Range("V14").Value = 'Södermanland'
$ws.Range("W14").Value = 'Mellösa'
$ws.Range("Y14").Value = '2015-03-02'
$ws.Range("Z14").Value = '00:00'
$ws.Range("AA14").Value = '2015-03-02'
$ws.Range("AB14").Value = '00:00'
$ws.Range("AC14").Value = 'Riklig förekomst på block.'
$ws.Range("AD14").Value = $False
$ws.Range("AE14").Value = $False
$ws.Range("AG14").Value = $False
$ws.Range("AT14").Value = ""
$ws.Range("AW14").Value = 'Bo Törnquist'
$ws.Range("AX14").Value = 'Bo Törnquist'
$ws.Range("AY14").Value = ""

# Row 15
$ws.Range("L15").ClearContents()
$ws.Range("A15").Value = 17177736
$ws.Range("B15").Value = 89652
$ws.Range("C15").Value = 'Ovaliderad'
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 73
$ws.Range("F15").Value = 'Veckticka'
$ws.Range("G15").Value = 'Flavidoporia pulvinascens'
$ws.Range("H15").Value = '(Pilát) Audet'
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = '1'
$ws.Range("J15").Value = 'mycel'
$ws.Range("K15").Value = ""
$ws.Range("P15").Value = 'S Kramnäs, Srm'
$ws.Range("Q15").Value = 596503.1774122283
$ws.Range("R15").Value = 6550127.803339794
$ws.Range("S15").Value = 5
$ws.Range("T15").Value = 'Södermanland'
$ws.Range("U15").Value = 'Flen'
$ws.Range("V15").Value = 'Södermanland'
$ws.Range("W15").Value = 'Mellösa'
$ws.Range("Y15").Value = '2015-03-02'
$ws.Range("Z15").Value = '00:00'
$ws.Range("AA15").Value = '2015-03-02'
$ws.Range("AB15").Value = '00:00'
$ws.Range("AC15").Value = 'Grov asplåga.'
$ws.Range("AD15").Value = $False
$ws.Range("AE15").Value = $False
$ws.Range("AF15").Value = ""
$ws.Range("AG15").Value = $False
$ws.Range("AT15").Value = ""
$ws.Range("AW15").Value = 'Bo Törnquist'
$ws.Range("AX15").Value = 'Bo Törnquist'
$ws.Range("AY15").Value = ""

# Row 16
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("A16").Value = 17177737
$ws.Range("B16").Value = 90138
$ws.Range("C16").Value = 'Ovaliderad'
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 366
$ws.Range("F16").Value = 'Kandelabersvamp'
$ws.Range("G16").Value = 'Artomyces pyxidatus'
$ws.Range("H16").Value = '(Pers.) Jülich'
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("P16").Value = 'S Kramnäs, Srm'
$ws.Range("Q16").Value = 596503.1774122283
$ws.Range("R16").Value = 6550127.803339794
$ws.Range("S16").Value = 5
$ws.Range("T16").Value = 'Södermanland'
$ws.Range("U16").Value = 'Flen'
$ws.Range("V16").Value = 'Södermanland'
$ws.Range("W16").Value = 'Mellösa'
$ws.Range("Y16").Value = '2015-03-02'
$ws.Range("Z16").Value = '00:00'
$ws.Range("AA16").Value = '2015-03-02'
$ws.Range("AB16").Value = '00:00'
$ws.Range("AC16").Value = 'Grov asplåga.'
$ws.Range("AD16").Value = $False
$ws.Range("AE16").Value = $False
$ws.Range("AF16").Value = ""
$ws.Range("AG16").Value = $False
$ws.Range("AT16").Value = ""
$ws.Range("AW16").Value = 'Bo Törnquist'
$ws.Range("AX16").Value = 'Bo Törnquist'
$ws.Range("AY16").Value = ""
